# إضافة حدث جديد في Card15
# 1) Fill previously-empty cells in row 15 with the literal text "nan"
#    (matching the convention already used throughout this sheet).
# 2) Append a new row 16 describing the new service event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# --- Step 1: create row 16 first, while row 15 still holds its original ----
# (blank) values. Copy row 15's formatting/cell-types down to row 16 so that
# the new row keeps the workbook's existing "everything is text" convention
# without Excel re-typing numeric-looking values (e.g. "15") as numbers.
$ws.Range("A15:R15").Copy($ws.Range("A16:R16"))

# Now overwrite row 16 with the actual new-event values.
$ws.Cells.Item(16, 12).Value = "8\3\2025"                              # L16 Date
$ws.Cells.Item(16, 14).Value = "تم تغيير الجرائد الخلفيه (5_8)"        # N16 Correction
$ws.Cells.Item(16, 15).Value = "الخبير"                                # O16 Serviced by

# --- Step 2: update row 15, filling its previously-empty cells with "nan" --
$row15NanCols = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 13, 16, 17, 18)  # B..K, M, P, Q, R
foreach ($col in $row15NanCols) {
    $ws.Cells.Item(15, $col).Value = "nan"
}

Write-Host "Row 16 added to Card15"
